# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46035

$ws.Range("B2").Value = 83.36
$ws.Range("C2").Value = 80.04000000000001
$ws.Range("D2").Value = 78.47
$ws.Range("E2").Value = 75.89
$ws.Range("F2").Value = 74.66
$ws.Range("G2").Value = 77.05
$ws.Range("H2").Value = 86.7
$ws.Range("I2").Value = 97.56
$ws.Range("J2").Value = 106.58
$ws.Range("K2").Value = 101.64
$ws.Range("L2").Value = 95.03
$ws.Range("M2").Value = 87.3
$ws.Range("N2").Value = 80.67
$ws.Range("O2").Value = 76.76000000000001
$ws.Range("P2").Value = 78.55
$ws.Range("Q2").Value = 87.67
$ws.Range("R2").Value = 92.29000000000001
$ws.Range("S2").Value = 102.76
$ws.Range("T2").Value = 114.19
$ws.Range("U2").Value = 116.75
$ws.Range("V2").Value = 140.65
$ws.Range("W2").Value = 135.56
$ws.Range("X2").Value = 117.08
$ws.Range("Y2").Value = 109.32
$ws.Range("Z2").Value = 95.69

$ws.Range("AB2").Value = 125.65
$ws.Range("AD2").Value = 138.11
$ws.Range("AF2").Value = 115.47

$ws.Range("AG2").Value = "0h-16h"
